$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 50 new rows starting at row 9, shifting existing rows 9+ down
$ws.Rows("9:58").Insert()

# Fill the newly inserted rows with new comment text
$ws.Range("A9").Value = 'I don’t think I would buy this.'
$ws.Range("A10").Value = 'Too many features can be confusing.'
$ws.Range("A11").Value = 'Not convinced it’s really useful.'
$ws.Range("A12").Value = 'Sounds like marketing hype.'
$ws.Range("A13").Value = 'I doubt it will give accurate advice.'
$ws.Range("A14").Value = 'Might be hard for older people.'
$ws.Range("A15").Value = 'I prefer simple fitness bands.'
$ws.Range("A16").Value = 'Probably needs too many subscriptions.'
$ws.Range("A17").Value = 'I don’t feel comfortable sharing health data.'
$ws.Range("A18").Value = 'This looks complicated to use.'
$ws.Range("A19").Value = 'Is there a demo video?'
$ws.Range("A20").Value = 'Can it be used while sleeping?'
$ws.Range("A21").Value = 'What is the return policy?'
$ws.Range("A22").Value = 'Is international shipping available?'
$ws.Range("A23").Value = 'Does it come with a warranty?'
$ws.Range("A24").Value = 'Can multiple users use one device?'
$ws.Range("A25").Value = 'How long does the battery last?'
$ws.Range("A26").Value = 'Does it support Android and iOS?'
$ws.Range("A27").Value = 'Is it available in black color?'
$ws.Range("A28").Value = 'What phones are compatible with it?'
$ws.Range("A29").Value = 'Looks like a thoughtful invention.'
$ws.Range("A30").Value = 'Seems well designed.'
$ws.Range("A31").Value = 'Pocket wellness coach — nice!'
$ws.Range("A32").Value = 'This could improve productivity.'
$ws.Range("A33").Value = 'I like the stress monitoring feature.'
$ws.Range("A34").Value = 'Great thinking behind this product.'
$ws.Range("A35").Value = 'Makes self-care easier.'
$ws.Range("A36").Value = 'This is such a smart innovation.'
$ws.Range("A37").Value = 'Looks very user-friendly.'
$ws.Range("A38").Value = 'Wellness tech done right.'
$ws.Range("A39").Value = 'Could be helpful for remote workers.'
$ws.Range("A40").Value = 'Love the personal assistant feel.'
$ws.Range("A41").Value = 'A smart solution for modern stress.'
$ws.Range("A42").Value = 'Helps people focus on their mental health.'
$ws.Range("A43").Value = 'Good product for personal wellness.'
$ws.Range("A44").Value = 'Feels modern and useful.'
$ws.Range("A45").Value = 'Would love to test this.'
$ws.Range("A46").Value = 'I like the sleek gadget idea.'
$ws.Range("A47").Value = 'This could reduce daily anxiety.'
$ws.Range("A48").Value = 'A nice blend of tech and health.'
$ws.Range("A49").Value = 'Sounds like a digital wellness buddy.'
$ws.Range("A50").Value = 'Smart reminders are a great touch.'
$ws.Range("A51").Value = 'I’d gift this to my parents.'
$ws.Range("A52").Value = 'Looks super handy to carry around.'
$ws.Range("A53").Value = 'Perfect for people working long hours.'
$ws.Range("A54").Value = 'I like the idea of mood monitoring.'
$ws.Range("A55").Value = 'This could really help people who forget self-care.'
$ws.Range("A56").Value = 'Wellness made simple — nice concept.'
$ws.Range("A57").Value = 'Love how compact yet powerful it sounds.'
$ws.Range("A58").Value = 'This gadget feels very premium and thoughtful.'

# Remove the two trailing rows that were dropped (now at 102 and 103)
$ws.Rows("102:103").Delete()

